# Update crypto price/volume snapshot (GitHub Actions scheduled refresh).
# Price values are textual (not numeric) in the source data; some new values
# look like plain numbers to Excel's auto-detection, so NumberFormat is
# forced to Text ("@") immediately before assignment to keep them as strings
# (matching the original inline-string cell content) instead of being
# silently coerced to floats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.824.57'
$ws.Range("E2").Value = '  +0.91%  '

$ws.Range("D3").Value = '1.807.63'
$ws.Range("E3").Value = '  +0.28%  '

$ws.Range("E4").Value = '  +0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.77'
$ws.Range("E5").Value = '  +2.73%  '

$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("E7").Value = '  +0.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.63'
$ws.Range("E8").Value = '  -6.47%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.310'
$ws.Range("E9").Value = '  +6.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0681'
$ws.Range("E10").Value = '  +2.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0997'
$ws.Range("E11").Value = '  +0.33%  '

$ws.Range("D12").Value = '2.071.36'
$ws.Range("E12").Value = '  +0.39%  '

$ws.Range("D13").Value = '1.793.13'
$ws.Range("E13").Value = '  -0.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.00'
$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.64'
$ws.Range("E15").Value = '  +5.40%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.653'
$ws.Range("E16").Value = '  +3.95%  '

$ws.Range("D17").Value = '34.813.05'
$ws.Range("E17").Value = '  +0.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.69'
$ws.Range("E18").Value = '  +2.01%  '

$ws.Range("D19").Value = '0.0₃0781'
$ws.Range("E19").Value = '  +1.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.61'
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.70'
$ws.Range("E21").Value = '  +4.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.63'
$ws.Range("E22").Value = '  +6.34%  '

$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("E24").Value = '  +5.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.89'
$ws.Range("E25").Value = '  +0.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.78'
$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.35'
$ws.Range("E27").Value = '  -0.38%  '

$ws.Range("E28").Value = '  -1.14%  '

$ws.Range("E29").Value = '  +28.83%  '

$ws.Range("E30").Value = '  +0.52%  '

$ws.Range("D31").Value = '3.339.43'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0547'
$ws.Range("E32").Value = '  +6.51%  '

$ws.Range("E33").Value = '  +1.65%  '

$ws.Range("E34").Value = '  +1.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.77'
$ws.Range("E35").Value = '  -0.81%  '

$ws.Range("E36").Value = '  +10.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '93.02'
$ws.Range("E37").Value = '  +6.11%  '

$ws.Range("E38").Value = '  +4.04%  '

$ws.Range("E39").Value = '  +1.91%  '

$ws.Range("D40").Value = '1.299.52'
$ws.Range("E40").Value = '  -1.60%  '

$ws.Range("E41").Value = '  +3.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.981'
$ws.Range("E42").Value = '  +4.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.74'
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("E44").Value = '  -1.39%  '

$ws.Range("E46").Value = '  -1.63%  '

$ws.Range("E47").Value = '  +7.55%  '

$ws.Range("E48").Value = '  -1.42%  '

$ws.Range("D49").Value = '1.985.61'
$ws.Range("E49").Value = '  +1.03%  '

$ws.Range("E50").Value = '  +0.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0643'
$ws.Range("E51").Value = '  +5.82%  '
